# Remove embedded line-breaks from a handful of vaccine name / note strings,
# turning multi-line cell text into single-line text (newline -> space).
#
# "Pediatric Influenza Vaccine " sheet
$wb = $excel.ActiveWorkbook
$wsPed = $wb.Worksheets.Item("Pediatric Influenza Vaccine ")
$wsPed.Range("B3").Value2  = "Fluzone Pediatric dose No Preservative"
$wsPed.Range("B6").Value2  = "Fluarix Preservative-Free"
$wsPed.Range("B9").Value2  = "FluMist No Preservative"
$wsPed.Range("B10").Value2 = "Afluria No Preservative"
$wsPed.Range("H10").Value2 = "Merck (CSL product)"

# "Adult Influenza Vaccine " sheet
$wsAdult = $wb.Worksheets.Item("Adult Influenza Vaccine ")
$wsAdult.Range("B5").Value2  = "Agriflu No Preservative"
$wsAdult.Range("B7").Value2  = "Fluvirin Preservative-free"
$wsAdult.Range("B8").Value2  = "Fluarix Preservative-free"
$wsAdult.Range("B10").Value2 = "Flumist No Preservative"
